# Update the "想去人数" (want-to-go count) values in column F, rows 2-6,
# on both the "展览" and "全部类型" worksheets (they contain the same data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2 = 373
    3 = 1263
    4 = 1542
    5 = 60
    6 = 6165
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
